# Apply the "more data, first round regression models" edit:
#  - Add 30 years of historical data (1950-1979) above the existing
#    1980-2019 data, extending the table to row 71.
#  - Update the 2018 yield value from 52.3 to 46.5.
#  - Set column B width and the active cell selection to match the target.
#
# NOTE: we intentionally avoid Rows.Insert()/EntireRow.Insert() here - that
# operation causes the runtime to rewrite every shifted cell's value as a
# full-precision double (e.g. 43.7 -> 43.700000000000003) even for cells
# that were never touched. Instead we just write the full, final data set
# directly into place with plain Value assignments, which round-trips
# cleanly and only changes precision for the literals we actually supply
# (matching how Excel itself stores those doubles, e.g. 34.2 ->
# 34.200000000000003, as already present in the original file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data set for A2:B71 (year, yield) -- years 1950 through 2019.
$allData = @(
    @(1950, 20.2),
    @(1951, 19.8),
    @(1952, 19.3),
    @(1953, 19.7),
    @(1954, 20),
    @(1955, 20),
    @(1956, 20.6),
    @(1957, 21.5),
    @(1958, 22.2),
    @(1959, 22.6),
    @(1960, 23.3),
    @(1961, 25.1),
    @(1962, 24),
    @(1963, 23),
    @(1964, 22.8),
    @(1965, 24.6),
    @(1966, 25),
    @(1967, 23.9),
    @(1968, 25.2),
    @(1969, 26.3),
    @(1970, 26.5),
    @(1971, 27.1),
    @(1972, 27.1),
    @(1973, 23.4),
    @(1974, 22.5),
    @(1975, 26.3),
    @(1976, 25.2),
    @(1977, 28.4),
    @(1978, 27.4),
    @(1979, 32.5),
    @(1980, 25.8),
    @(1981, 28.3),
    @(1982, 27.5),
    @(1983, 25.7),
    @(1984, 26.9),
    @(1985, 33.6),
    @(1986, 34.2),
    @(1987, 30),
    @(1988, 32.5),
    @(1989, 32.5),
    @(1990, 34),
    @(1991, 34.7),
    @(1992, 37.5),
    @(1993, 32.5),
    @(1994, 42.3),
    @(1995, 35.6),
    @(1996, 37.5),
    @(1997, 49),
    @(1998, 49),
    @(1999, 36.2),
    @(2000, 38),
    @(2001, 39.5),
    @(2002, 48),
    @(2003, 34),
    @(2004, 43.7),
    @(2005, 44.2),
    @(2006, 43),
    @(2007, 42),
    @(2008, 39.5),
    @(2009, 44.6),
    @(2010, 44),
    @(2011, 42),
    @(2012, 40),
    @(2013, 44.6),
    @(2014, 47.8),
    @(2015, 48),
    @(2016, 52.6),
    @(2017, 49.5),
    @(2018, 46.5),
    @(2019, 46.9)
)

for ($i = 0; $i -lt $allData.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $allData[$i][0]
    $ws.Cells.Item($row, 2).Value = $allData[$i][1]
}

# Column B width adjustment to fit the longer header/values. (The headless
# COM runtime quantizes ColumnWidth to 1/6-character increments, so 23.33
# is the closest input that reproduces the target stored width of
# ~24.1640625 characters.)
$ws.Columns.Item(2).ColumnWidth = 23.33

# Update the selected/active cell to match the saved view state.
$ws.Range("E18").Select()
